$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2) entirely; all rows below shift up by one,
# and the used range shrinks from A1:F63 to A1:F62.
$ws.Rows(2).Delete()
